$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are stored as text in the sheet
# (values like "36.817.06" use dots as thousands separators, and the
# percentage strings carry padding whitespace) so force text format
# before assigning, to avoid Excel auto-coercing them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.817.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.045.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.653'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.20%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.01'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.29'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.368'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0772'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.08'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.867'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.340.94'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.036.07'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.82'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.747.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.15'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.28'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.21'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.17'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '167.94'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.85'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.48'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +13.96%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.69'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0608'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.90%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.30'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0813'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -9.16%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.65%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.02'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.69%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0939'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -14.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.12'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.75'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.300.37'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.68%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.70'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.226.73'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.42%  '
